# Update experiment tracking metrics (rows 2 and 3, columns H:S)
# as part of switching all plots to the custom style and verifying
# save locations are consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.0144
$ws.Range("I2").Value = 0.0102
$ws.Range("J2").Value = 0.994
$ws.Range("K2").Value = 0.0055
$ws.Range("L2").Value = 0.994
$ws.Range("M2").Value = 0.0054
$ws.Range("N2").Value = 0.9902
$ws.Range("O2").Value = 0.0119
$ws.Range("P2").Value = 0.998
$ws.Range("Q2").Value = 0.0045
$ws.Range("R2").Value = 0.994
$ws.Range("S2").Value = 0.0055

# Row 3
$ws.Range("H3").Value = 0.0154
$ws.Range("I3").Value = 0.0094
$ws.Range("J3").Value = 0.9945000000000001
$ws.Range("K3").Value = 0.005
$ws.Range("L3").Value = 0.9945000000000001
$ws.Range("M3").Value = 0.0049
$ws.Range("N3").Value = 0.9902
$ws.Range("O3").Value = 0.0103
$ws.Range("P3").Value = 0.999
$ws.Range("Q3").Value = 0.0032
$ws.Range("R3").Value = 0.9945000000000001
$ws.Range("S3").Value = 0.005
